$d = $word.ActiveDocument

# Locate the paragraph that ends the German judgement-example bullet
# ("... B: Nein, nach (Afrika?).") so the new bullet can be inserted
# right after it, before the "Methods:" paragraph.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Nein, nach (Afrika?)*") {
        $anchor = $p
    }
}

# A collapsed range positioned right before the anchor paragraph's
# paragraph mark (End - 1). Calling InsertXML there inserts a brand
# new sibling <w:p> right after the anchor paragraph, with none of the
# anchor's run/paragraph formatting (e.g. lang="de-DE") carried over -
# exactly matching a freshly authored paragraph in the target markup.
$insertPoint = $d.Range($anchor.Range.End - 1, $anchor.Range.End - 1)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Perhaps, other stimuli that a</w:t></w:r><w:r><w:t xml:space="preserve">re less difficult to translate </w:t></w:r><w:r><w:t xml:space="preserve">(no issues with verbals froms and prepositions) </w:t></w:r><w:r><w:t>and still show contrastive focus, e. g. A: Ich habe geh' + [char]0x00F6 + 'rt, Sandy spielt Fu' + [char]0x00DF + 'ball. B: Nein, Tennis (</w:t></w:r><w:r><w:t>cf</w:t></w:r><w:r><w:t>. Konietzko &amp; Winkler</w:t></w:r><w:r><w:t xml:space="preserve"> 2010: 1437</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>'

$null = $insertPoint.InsertXML($newParaXml)

Write-Host "Inserted new bullet paragraph."
